$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) New empty paragraph right after the existing final paragraph
#    (spacing: after=0, line=360/auto -- matches the surrounding body
#    text paragraph already in the document).
# ------------------------------------------------------------------
$endPos1 = $d.Content.End - 1
$insertPoint1 = $d.Range($endPos1, $endPos1)
$insertPoint1.InsertParagraphAfter()

$emptyPara = $d.Paragraphs.Last
$emptyPara.Format.SpaceAfter = 0
$emptyPara.Format.LineSpacingRule = 1

# ------------------------------------------------------------------
# 2) New paragraph containing the bold "Feedback" heading.
# ------------------------------------------------------------------
$endPos2 = $d.Content.End - 1
$insertPoint2 = $d.Range($endPos2, $endPos2)
$insertPoint2.InsertParagraphAfter()

$feedbackPara = $d.Paragraphs.Last
$feedbackPara.Format.SpaceAfter = 0
$feedbackPara.Format.LineSpacingRule = 1

$feedbackTextPos = $d.Content.End - 1
$feedbackTextRange = $d.Range($feedbackTextPos, $feedbackTextPos)
$feedbackTextRange.Text = "Feedback"

# ------------------------------------------------------------------
# 3) New paragraph with Brother Birch's comment (typed before the
#    bold formatting is applied to the Feedback paragraph so this
#    new paragraph does not inherit the bold formatting).
# ------------------------------------------------------------------
$endPos3 = $d.Content.End - 1
$insertPoint3 = $d.Range($endPos3, $endPos3)
$insertPoint3.InsertParagraphAfter()

$commentPara = $d.Paragraphs.Last
$commentPara.Format.SpaceAfter = 0
$commentPara.Format.LineSpacingRule = 1

$commentTextPos = $d.Content.End - 1
$commentTextRange = $d.Range($commentTextPos, $commentTextPos)
$commentTextRange.Text = "Brother Birch: It looks great, honestly wouldn’t change anything."

# ------------------------------------------------------------------
# Now apply bold formatting to the "Feedback" paragraph (including
# its paragraph mark) now that it is safely bounded by paragraph 3.
# ------------------------------------------------------------------
$feedbackPara.Range.Font.Bold = $true
$feedbackPara.Range.Font.BoldBi = $true
